# "move to new server" - add newly listed products to the price sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Бельгийский шоколад 3 плитки
$ws.Range("A3").Value = "Бельгийский шоколад 3 плитки "
$ws.Range("B3").Value = 5500
$ws.Range("C3").Value = 50

# Row 4: 0.33 унции бельгийских камней
$ws.Range("A4").Value = "0.33 унции бельгийских камней"
$ws.Range("B4").Value = 2200
$ws.Range("C4").Value = 52

# Row 5: 0.5 унции бельгийских камней
$ws.Range("A5").Value = "0.5 унции бельгийских камней"
$ws.Range("B5").Value = 3000
$ws.Range("C5").Value = 51

# Leave the selection on B5, matching the saved view state
[void]$ws.Range("B5").Select()
